$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Capture the *original* column S ("comms_internal") formatting into
#    the new column U before column S gets overwritten with the values
#    that currently live in column M ("comms"). Data rows only (2..131);
#    row 1 (header) is handled separately below.
# ---------------------------------------------------------------------
$ws.Range("S2:S131").Copy()
$ws.Range("U2:U131").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 2. Move column M ("comms") into column S ("comms_internal"):
#    - S adopts M's formatting (style)
#    - S adopts M's values
#    - M is cleared (but keeps its original formatting)
# ---------------------------------------------------------------------
$ws.Range("M2:M131").Copy()
$ws.Range("S2:S131").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("M2:M131").Copy()
$ws.Range("S2:S131").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("M2:M131").ClearContents()

# ---------------------------------------------------------------------
# 3. Give the new column T ("verb_original_orthography") data cells the
#    same formatting as column P (plain style used by the neighbouring
#    empty columns).
# ---------------------------------------------------------------------
$ws.Range("P2:P131").Copy()
$ws.Range("T2:T131").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 4. Header row: add the two new header cells, matching the formatting
#    used by the other rotated-text headers (e.g. D1).
# ---------------------------------------------------------------------
$ws.Range("D1").Copy()
$ws.Range("T1:U1").PasteSpecial(-4122)     # xlPasteFormats
$ws.Range("T1").Value = "verb_original_orthography"
$ws.Range("U1").Value = "sentence_original_orthography"

# ---------------------------------------------------------------------
# 5. Column widths for the two new columns, matching R (for T) and
#    S (for U).
# ---------------------------------------------------------------------
$ws.Columns.Item(20).ColumnWidth = 27.75
$ws.Columns.Item(21).ColumnWidth = 112.9167

# ---------------------------------------------------------------------
# 6. Refresh the AutoFilter so it spans the new columns.
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:U131").AutoFilter()

# ---------------------------------------------------------------------
# 7. Update the (hidden) _FilterDatabase defined name to match.
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Yakut!_FilterDatabase") {
        $n.RefersTo = "=Yakut!`$A`$1:`$U`$131"
    }
}

# ---------------------------------------------------------------------
# 8. Restore view state (scroll position / selection) as closely as
#    possible.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M9").Select()
